# Update existing "ODI Batting" sheet: rename MATCH_CARD_LINK -> MATCH_CODE
# and replace the full scorecard URLs with the bare match-code numbers.
$wb = $excel.ActiveWorkbook
$odi = $wb.Worksheets.Item("ODI Batting")

$odi.Range("D1").Value = "MATCH_CODE"

$odi.Range("D2").NumberFormat = "@"
$odi.Range("D2").Value = "4727"

$odi.Range("D3").NumberFormat = "@"
$odi.Range("D3").Value = "4731"

# Insert a brand-new "Player Info" sheet ahead of "ODI Batting" with the
# player's basic info (scraped alongside the batting performance data).
$info = $wb.Worksheets.Add()
$info.Name = "Player Info"

$info.Range("A1").Value = "ID"
$info.Range("B1").Value = "NAME"
$info.Range("C1").Value = "BATTING_HAND"
$info.Range("D1").Value = "BOWL_STYLE"

# Match the bold/centered/bordered header look used on the "ODI Batting" sheet.
$hdr = $info.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

$info.Range("A2").NumberFormat = "@"
$info.Range("A2").Value = "6552"
$info.Range("B2").Value = "Ryan David Rickelton"
$info.Range("C2").Value = "Left Handed"
$info.Range("D2").Value = "Does Not Bowl | Unknown"
